# The source workbook tracks monthly ISM Manufacturing PMI values in
# columns A (date, as a serial number) / B (index value). This adds the
# four newest months (Sep-Dec 2024) to the bottom of the existing table
# on Sheet1, continuing the same format as the rows immediately above.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert new rows one at a time, copying formatting from the row above
# (xlFormatFromLeftOrAbove) so the date/value styles (s="3"/s="1") carry
# forward exactly like the preceding rows, instead of falling back to
# generic/default formatting.
$ws.Rows.Item(368).Insert(-4121, -4163)
$ws.Rows.Item(369).Insert(-4121, -4163)
$ws.Rows.Item(370).Insert(-4121, -4163)
$ws.Rows.Item(371).Insert(-4121, -4163)

$ws.Cells.Item(368, 1).Value = 45565
$ws.Cells.Item(368, 2).Value = 48.3

$ws.Cells.Item(369, 1).Value = 45596
$ws.Cells.Item(369, 2).Value = 54.8

$ws.Cells.Item(370, 1).Value = 45626
$ws.Cells.Item(370, 2).Value = 50.3

$ws.Cells.Item(371, 1).Value = 45657
$ws.Cells.Item(371, 2).Value = 52.5

# Column A widened slightly to keep fitting the date strings.
$ws.Columns.Item(1).ColumnWidth = 9.619791666666666

# Move the selection down to the newly added last cell and scroll the
# window so the new rows are visible, matching the refreshed view.
$ws.Range("B371").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 357
